# Auto-generated edit script: update crypto price/volume table values
# per commit "Updated cryptos list on Fri Jun 30 10:51:46 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.819.46"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.78"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.25"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4800"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2957"
$ws.Range("E8").Value = "  +7.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06646"
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.66"
$ws.Range("E10").Value = "  +6.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "100.37"
$ws.Range("E11").Value = "  +18.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.886.86"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07585"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.172"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6638"
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "305.54"
$ws.Range("E16").Value = "  +27.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.788.38"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.20"
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007631"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.123.48"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9980"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.172"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.230"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.339"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.88"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.33"
$ws.Range("E27").Value = "  +12.09%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1161"
$ws.Range("E28").Value = "  +13.03%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.959"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.196"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.007"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05098"
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7470"
$ws.Range("E34").Value = "  +5.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01978"
$ws.Range("E37").Value = "  +3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.065"
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8985"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.33"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4217"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.659"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.97"
$ws.Range("E45").Value = "  +10.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.394"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.091"
$ws.Range("E47").Value = "  +5.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1235"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.08"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05647"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.401"
$ws.Range("E51").Value = "  +1.65%  "
